$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 24, shifting the old row 24 (week of 44701) down to row 25.
$ws.Rows.Item(24).Insert()

# Row 22: update with the new latest week's data (44798)
$ws.Range("D22").Value = 44798
$ws.Range("J22").Value = 200
$ws.Range("L22").Value = 8500
$ws.Range("M22").Value = 8250
$ws.Range("P22").Value = 2750

# Row 23: shift to the previous week's data (44771)
$ws.Range("D23").Value = 44771
$ws.Range("J23").Value = 150

# Row 24 (new row): fill with the week of 44782
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C24").Value = "Ñuble"
$ws.Range("D24").Value = 44782
$ws.Range("E24").Value = 16
$ws.Range("F24").Value = 100112037
$ws.Range("G24").Value = "Cebollín"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 8000
$ws.Range("L24").Value = 8000
$ws.Range("M24").Value = 8000
$ws.Range("N24").Value = "$/docena de atados"
$ws.Range("O24").Value = "Provincia de Diguillín"
$ws.Range("P24").Value = 2667
$ws.Range("Q24").Value = 3
$ws.Range("R24").Value = "Hortaliza"
